$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column I
$ws.Range("I1").Value = "Other found locations"

# Column I values for rows 2-7
$ws.Range("I2").Value = ""
$ws.Range("I3").Value = "_PMC"
$ws.Range("I4").Value = "_PMC"
$ws.Range("I5").Value = "_PMC"
$ws.Range("I6").Value = ""
$ws.Range("I7").Value = "_PMC"

# Updated Authors (column E) values for rows 3, 4, 5, 7
$ws.Range("E3").Value = "[Kenneth A.%Egol%NULL%0,    Sanjit R.%Konda%NULL%1,    Mackenzie L.%Bird%NULL%1,    Nicket%Dedhia%NULL%1,    Emma K.%Landes%NULL%1,    Rachel A.%Ranson%NULL%1,    Sara J.%Solasz%NULL%1,    Vinay K.%Aggarwal%NULL%1,    Joseph A.%Bosco%NULL%1,    David L.%Furgiuele%NULL%1,    Abhishek%Ganta%NULL%1,    Jason%Gould%NULL%1,    Thomas R.%Lyon%NULL%1,    Toni M.%McLaurin%NULL%1,    Nirmal C.%Tejwani%NULL%1,    Joseph D.%Zuckerman%NULL%1,    Philipp%Leucht%NULL%1]"
$ws.Range("E4").Value = "[Drake G.%LeBrun%NULL%0,    Maxwell A.%Konnaris%NULL%1,    Gregory C.%Ghahramani%NULL%1,    Ajay%Premkumar%NULL%1,    Chris J.%DeFrancesco%NULL%1,    Jordan A.%Gruskay%NULL%1,    Aleksey%Dvorzhinskiy%NULL%1,    Milan S.%Sandhu%NULL%1,    Elan M.%Goldwyn%NULL%1,    Christopher L.%Mendias%NULL%1,    William M.%Ricci%NULL%1]"
$ws.Range("E5").Value = "[Amit%Thakrar%NULL%1,    Karen%Chui%NULL%1,    Akhil%Kapoor%NULL%1,    John%Hambidge%NULL%1]"
$ws.Range("E7").Value = "[Karen%Chui%NULL%1,    Amit%Thakrar%NULL%1,    Shivakumar%Shankar%NULL%1]"
